$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 12 (shifts existing rows 12-46 down to 13-47)
$ws.Rows.Item(12).Insert()

# Populate the new row with the "hyphen" library entry
$ws.Range("B12").Value = "hyphen"
$ws.Hyperlinks.Add($ws.Range("C12"), "https://www.npmjs.com/package/hyphen")
$ws.Range("D12").Value = "This is a text hyphenation library, based on Franklin M. Liang's hyphenation algorithm. In core of the algorithm lies a set of hyphenation patterns. They are extracted from hand-hyphenated dictionaries. Patterns for this library were taken from ctan.org and ported to Javascript."

# Re-apply the named Hyperlink style so the cell reuses the existing style (not a new duplicate)
$ws.Range("C12").Style = "Hyperlink"

# Update the active selection to match the edited workbook state
$ws.Range("C9").Select()
